# Apply daily spot price update (row 2) as described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45940
$ws.Range("B2").Value = 89.02
$ws.Range("C2").Value = 85.02
$ws.Range("D2").Value = 79.98999999999999
$ws.Range("E2").Value = 72.95999999999999
$ws.Range("F2").Value = 69.75
$ws.Range("G2").Value = 78.87
$ws.Range("H2").Value = 94.98
$ws.Range("I2").Value = 100.54
$ws.Range("J2").Value = 107.24
$ws.Range("K2").Value = 93.84
$ws.Range("L2").Value = 48.55
$ws.Range("M2").Value = 42.92
$ws.Range("N2").Value = 43.91
$ws.Range("O2").Value = 41.08
$ws.Range("P2").Value = 33.23
$ws.Range("Q2").Value = 45.71
$ws.Range("R2").Value = 52.19
$ws.Range("S2").Value = 79.43000000000001
$ws.Range("T2").Value = 96.13
$ws.Range("U2").Value = 115.01
$ws.Range("V2").Value = 117.47
$ws.Range("W2").Value = 108.1
$ws.Range("X2").Value = 103.43
$ws.Range("Y2").Value = 101.96
$ws.Range("Z2").Value = 79.22

# AA2 unchanged ("20h-24h")
$ws.Range("AB2").Value = 107.74

# AC2 unchanged ("20h-22h")
$ws.Range("AD2").Value = 112.78

$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 105.57
$ws.Range("AG2").Value = "3h-16h"
